$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '93.545.62'
$ws.Range("E2").Value = '  +3.72%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.130.43'
$ws.Range("E3").Value = '  +0.59%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.83'
$ws.Range("E5").Value = '  +2.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.54'
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("E7").Value = '  +1.44%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.410'
$ws.Range("E8").Value = '  +11.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.127.88'
$ws.Range("E10").Value = '  +30.81%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.745'
$ws.Range("E11").Value = '  +1.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.202'
$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("E13").Value = '  +4.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.96'
$ws.Range("E14").Value = '  +0.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.304.40'
$ws.Range("E15").Value = '  +3.42%  '

$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.728.12'
$ws.Range("E17").Value = '  +0.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.103.27'
$ws.Range("E18").Value = '  -1.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.80'
$ws.Range("E19").Value = '  +4.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.94'
$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("B21").Value = 'PEPE'
$ws.Range("C21").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000209'
$ws.Range("E21").Value = '  +4.27%  '

$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.87'
$ws.Range("E22").Value = '  +0.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '451.86'
$ws.Range("E23").Value = '  +3.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.44'
$ws.Range("E24").Value = '  +5.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.86'
$ws.Range("E25").Value = '  +2.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.97'
$ws.Range("E26").Value = '  +7.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.90'
$ws.Range("E27").Value = '  +1.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.302.47'
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("E30").Value = '  +10.97%  '

$ws.Range("E31").Value = '  +1.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.229'
$ws.Range("E32").Value = '  +0.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.32'
$ws.Range("E33").Value = '  +1.11%  '

$ws.Range("E34").Value = '  +8.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.12'
$ws.Range("E35").Value = '  +6.86%  '

$ws.Range("E36").Value = '  -4.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.43'
$ws.Range("E37").Value = '  +1.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.92'
$ws.Range("E38").Value = '  +0.24%  '

$ws.Range("B39").Value = 'MantraDAO'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.88'
$ws.Range("E39").Value = '  +4.40%  '

$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '485.62'
$ws.Range("E40").Value = '  -2.89%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.31'
$ws.Range("E41").Value = '  -2.13%  '

$ws.Range("E42").Value = '  +1.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.438'
$ws.Range("E43").Value = '  -1.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.10'
$ws.Range("E44").Value = '  +4.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '162.92'
$ws.Range("E46").Value = '  +2.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.95'
$ws.Range("E47").Value = '  +2.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.697'
$ws.Range("E48").Value = '  -2.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.40'
$ws.Range("E49").Value = '  +3.36%  '

$ws.Range("E50").Value = '  +5.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.49'
$ws.Range("E51").Value = '  +2.29%  '
